# MEGATAB_EMPREEND_JUN2025vlight.xlsx — "Add files via upload"
#
# The COORDENADA(DEC) column (E) on Planilha1 was re-formatted from
# "lat, lng" (e.g. "-3.891234, -38.455678") to a "lng,lat,0" KML/GIS style
# triplet (e.g. "-38.455678,-3.891234,0"). Apply that transform to every
# populated coordinate cell, then restore the sheet's default view
# (scrolled/selected at A1) and size the now-adjacent column P the way the
# refreshed sheet expects.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Reformat COORDENADA(DEC) values in column E (rows 2-23) -----------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $old = $cell.Text
    if ([string]::IsNullOrEmpty($old)) { continue }

    $parts = $old.Split(",")
    if ($parts.Length -ne 2) { continue }

    $lat = $parts[0].Trim()
    $lng = $parts[1].Trim()
    $cell.Value = "$lng,$lat,0"
}

# --- Reset the sheet view back to the top-left / A1 ---------------------
$ws.Activate()
$ws.Range("A1").Select()

# --- Give the new neighboring column (P) a sensible, fitted width -------
$ws.Columns.Item(16).ColumnWidth = 19.86

$wb.Save()
